$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------
# Sheet "Raw Material" (sheet1): rename table header cells
#   Kode  -> Kode Material
#   Specs -> Spesifikasi
# ---------------------------------------------------------------
$ws1.Cells.Item(1,1).Value2 = "Kode Material"
$ws1.Cells.Item(1,3).Value2 = "Spesifikasi"

# ---------------------------------------------------------------
# Sheet "Progress" (sheet2): add a 3rd "Deskripsi" column to Table2
# and fill it in with the process descriptions
# ---------------------------------------------------------------
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListColumns.Add() | Out-Null

# bring over the header formatting (fill/bold) from column B's header
$ws2.Cells.Item(1,2).Copy() | Out-Null
$ws2.Cells.Item(1,3).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws2.Cells.Item(1,3).Value2 = "Deskripsi"
$ws2.Cells.Item(2,3).Value2 = "Pemotongan Bahan"
$ws2.Cells.Item(3,3).Value2 = "Penganuan Bahan"
$ws2.Cells.Item(4,3).Value2 = "Pengecekan Kualitas Pemotongan dan Penganuan Barang"

$ws2.Columns.Item(3).AutoFit() | Out-Null

$ws2.PageSetup.Orientation = 1   # xlPortrait

# ---------------------------------------------------------------
# Selections left by the author when the workbook was last saved
# ---------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B13").Select() | Out-Null

$ws2.Activate()
$ws2.Range("D11").Select() | Out-Null
